# Weekly update: insert a new Ciboulette price record for
# "Vega Modelo de Temuco" right after the existing row for 2021-10-05
# (row 238), which pushes all subsequent rows down by one and appends
# the former last row (old row 255) as the new row 256.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 238 (existing row 238 and below
# shift down to 239 and below).
$ws.Rows.Item(238).Insert()

# Populate the newly inserted row 238 with the new record's data.
$ws.Cells.Item(238, 1).Value = 10
$ws.Cells.Item(238, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(238, 3).Value = "La Araucanía"
$ws.Cells.Item(238, 4).Value = 44746
$ws.Cells.Item(238, 5).Value = 9
$ws.Cells.Item(238, 6).Value = 100112039
$ws.Cells.Item(238, 7).Value = "Ciboulette"
$ws.Cells.Item(238, 8).Value = "Sin especificar"
$ws.Cells.Item(238, 9).Value = "Primera"
$ws.Cells.Item(238, 10).Value = 55
$ws.Cells.Item(238, 11).Value = 9000
$ws.Cells.Item(238, 12).Value = 9000
$ws.Cells.Item(238, 13).Value = 9000
$ws.Cells.Item(238, 14).Value = "`$/docena de atados"
$ws.Cells.Item(238, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(238, 16).Value = 3000
$ws.Cells.Item(238, 17).Value = 3
$ws.Cells.Item(238, 18).Value = "Hortaliza"
